$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking value to be stored as TEXT (matches source
# inlineStr cells) without leaving a stray NumberFormat on the cell -
# apostrophe-prefix forces text, then resetting the Style back to Normal
# drops the quotePrefix style bump so the cell keeps the default style.
function Set-TextValue($cell, $text) {
    $cell.Formula = "'" + $text
    $cell.Style = "Normal"
}

# Row 2
$ws.Range("B2").Value = 44030
$ws.Range("C2").Value = 27987
$ws.Range("D2").Value = 749
$ws.Range("E2").Value = 3222
$ws.Range("G2").Value = 11.51
$ws.Range("H2").Value = 12.68

# Row 3
$ws.Range("B3").Value = 44030
$ws.Range("C3").Value = 28633
$ws.Range("D3").Value = 251

# Row 4
$ws.Range("B4").Value = 44030
Set-TextValue $ws.Range("C4") "217895"
Set-TextValue $ws.Range("D4") "18771"
$ws.Range("E4").Value = 33686
$ws.Range("F4").Value = 5251
$ws.Range("G4").Value = 30.1
$ws.Range("K4").Value = 111896
$ws.Range("L4").Value = 17221

# Row 7
$ws.Range("B7").Value = 44030
$ws.Range("C7").Value = 76336
$ws.Range("D7").Value = 838
$ws.Range("E7").Value = 15111
$ws.Range("F7").Value = 300
$ws.Range("G7").Value = 19.8
$ws.Range("H7").Value = 35.8

# Row 8
$ws.Range("B8").Value = 44030
Set-TextValue $ws.Range("C8") "33332"
Set-TextValue $ws.Range("D8") "243"
Set-TextValue $ws.Range("E8") "830"
Set-TextValue $ws.Range("F8") "3"
$ws.Range("H8").Value = 1.23

# Row 9
$ws.Range("B9").Value = 44030
$ws.Range("C9").Value = 22184
$ws.Range("D9").Value = 667
$ws.Range("E9").Value = 2343
$ws.Range("F9").Value = 26
$ws.Range("G9").Value = 15.49
$ws.Range("H9").Value = 4.2
$ws.Range("K9").Value = 15127
$ws.Range("L9").Value = 622

# Row 10
$ws.Range("C10").Value = 32533
$ws.Range("D10").Value = 357
$ws.Range("E10").Value = 6918
$ws.Range("F10").Value = 93
$ws.Range("H10").Value = 26.2
$ws.Range("K10").Value = 28197
$ws.Range("L10").Value = 355

# Row 11
$ws.Range("B11").Value = 44030
$ws.Range("C11").Value = 23114
$ws.Range("D11").Value = 478
$ws.Range("E11").Value = 847
$ws.Range("G11").Value = 4.68
$ws.Range("H11").Value = 3.9
$ws.Range("K11").Value = 18089
$ws.Range("L11").Value = 462

# Row 13
$ws.Range("B13").Value = 44030
$ws.Range("C13").Value = 16736
$ws.Range("D13").Value = 569
$ws.Range("E13").Value = 303

# Row 16
$ws.Range("B16").Value = 44030
$ws.Range("C16").Value = 64180
$ws.Range("D16").Value = 1253
$ws.Range("E16").Value = 19494
$ws.Range("F16").Value = 536
$ws.Range("G16").Value = 43.92
$ws.Range("H16").Value = 44.67
$ws.Range("K16").Value = 44381
$ws.Range("L16").Value = 1200

# Row 17
$ws.Range("B17").Value = 44029
$ws.Range("C17").Value = 153041
$ws.Range("D17").Value = 4084
$ws.Range("E17").Value = 4094
$ws.Range("F17").Value = 409
$ws.Range("G17").Value = 4.69
$ws.Range("H17").Value = 10.76
$ws.Range("K17").Value = 87304
$ws.Range("L17").Value = 3801

# Row 19
$ws.Range("B19").Value = 44029
$ws.Range("C19").Value = 41846
$ws.Range("D19").Value = 1346
$ws.Range("E19").Value = 19138
$ws.Range("F19").Value = 671
$ws.Range("G19").Value = 45.73

# Row 22
$ws.Range("B22").Value = 44030
$ws.Range("C22").Value = 2471
$ws.Range("E22").Value = 13
$ws.Range("G22").Value = 0.53

# Row 26
$ws.Range("B26").Value = 44030
$ws.Range("C26").Value = 39788
$ws.Range("D26").Value = 1752
$ws.Range("E26").Value = 1981
$ws.Range("G26").Value = 6.25
$ws.Range("H26").Value = 6.96
$ws.Range("K26").Value = 31715
$ws.Range("L26").Value = 1695

# Row 27
$ws.Range("B27").Value = 44030
$ws.Range("C27").Value = 22481

# Row 28
$ws.Range("B28").Value = 44030
$ws.Range("C28").Value = 73098
$ws.Range("D28").Value = 6039
$ws.Range("E28").Value = 21215
$ws.Range("F28").Value = 2406
$ws.Range("G28").Value = 29.02
$ws.Range("H28").Value = 39.84

# Row 29
$ws.Range("B29").Value = 44029
$ws.Range("C29").Value = 375363
$ws.Range("D29").Value = 7595
$ws.Range("E29").Value = 10432
$ws.Range("F29").Value = 641
$ws.Range("G29").Value = 4.32
$ws.Range("H29").Value = 8.69
$ws.Range("K29").Value = 241390
$ws.Range("L29").Value = 7376

# Row 30
$ws.Range("B30").Value = 44030
$ws.Range("C30").Value = 55654
$ws.Range("D30").Value = 2627
$ws.Range("E30").Value = 6493
$ws.Range("G30").Value = 11.67
$ws.Range("H30").Value = 14.2

# Row 31
$ws.Range("B31").Value = 44030
$ws.Range("C31").Value = 1795
$ws.Range("D31").Value = 18
$ws.Range("E31").Value = 42
$ws.Range("G31").Value = 1.4
$ws.Range("K31").Value = 3010
$ws.Range("L31").Value = 36

# Row 32
$ws.Range("B32").Value = 44030
$ws.Range("C32").Value = 41485
$ws.Range("D32").Value = 843
$ws.Range("E32").Value = 6721
$ws.Range("G32").Value = 17.93
$ws.Range("H32").Value = 23.64
$ws.Range("K32").Value = 37485
$ws.Range("L32").Value = 829

# Row 33
$ws.Range("B33").Value = 44030
$ws.Range("C33").Value = 139872
$ws.Range("D33").Value = 3168
$ws.Range("E33").Value = 36504
$ws.Range("F33").Value = 1470
$ws.Range("G33").Value = 26.1
$ws.Range("H33").Value = 46.4

# Row 34
$ws.Range("B34").Value = 44030
$ws.Range("C34").Value = 46026
$ws.Range("D34").Value = 1444
$ws.Range("E34").Value = 1690
$ws.Range("G34").Value = 5.42
$ws.Range("H34").Value = 3.46
$ws.Range("K34").Value = 31176
$ws.Range("L34").Value = 1386

# Row 38
$ws.Range("C38").Value = 38197
$ws.Range("D38").Value = 789
$ws.Range("E38").Value = 3188
$ws.Range("G38").Value = 8.35
$ws.Range("H38").Value = 4.82

# Row 39
$ws.Range("B39").Value = 44030
$ws.Range("C39").Value = 97958
$ws.Range("D39").Value = 1629
$ws.Range("E39").Value = 15936
$ws.Range("F39").Value = 516
$ws.Range("G39").Value = 23.89
$ws.Range("H39").Value = 32.8
$ws.Range("K39").Value = 66695
$ws.Range("L39").Value = 1573

# Row 40
$ws.Range("B40").Value = 44030
$ws.Range("C40").Value = 160610
$ws.Range("D40").Value = 7290
$ws.Range("E40").Value = 27009
$ws.Range("F40").Value = 2011
$ws.Range("G40").Value = 16.82

# Row 41
$ws.Range("B41").Value = 44030
$ws.Range("C41").Value = 14302
$ws.Range("E41").Value = 177
$ws.Range("G41").Value = 1.24

# Row 42
$ws.Range("B42").Value = 44030
$ws.Range("C42").Value = 45470
$ws.Range("D42").Value = 1538
$ws.Range("E42").Value = 9200
$ws.Range("F42").Value = 151
$ws.Range("G42").Value = 20.23
$ws.Range("H42").Value = 9.82

# Row 43
$ws.Range("B43").Value = 44030
$ws.Range("C43").Value = 113238
$ws.Range("D43").Value = 8419
$ws.Range("E43").Value = 10642
$ws.Range("F43").Value = 691

# Row 45
$ws.Range("B45").Value = 44030
$ws.Range("C45").Value = 32246
$ws.Range("D45").Value = 1130
$ws.Range("E45").Value = 7584
$ws.Range("F45").Value = 380
$ws.Range("G45").Value = 32.77
$ws.Range("H45").Value = 36.09
$ws.Range("K45").Value = 23144
$ws.Range("L45").Value = 1053
